$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 467.21738
$ws.Range("I58").Value = 240.4375
$ws.Range("K58").Value = 721.3125
$ws.Range("M58").Value = -571.3125

$ws.Range("H62").Value = 2550
$ws.Range("I62").Value = 1600
$ws.Range("K62").Value = 1600
$ws.Range("M62").Value = -976

$ws.Range("H65").Value = 2550
$ws.Range("I65").Value = 1600
$ws.Range("K65").Value = 8000
$ws.Range("M65").Value = -4880

$ws.Range("H87").Value = 15460.454
$ws.Range("J87").Value = 15460.454
$ws.Range("L87").Value = 15460.454
$ws.Range("N87").Value = -17956.454

$ws.Range("H90").Value = 15460.454
$ws.Range("J90").Value = 15460.454
$ws.Range("L90").Value = 46381.362
$ws.Range("N90").Value = -58861.362

$ws.Range("H103").Value = 343.2353
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 391.875
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 1175.625
$ws.Range("M103").Value = -314
$ws.Range("N103").Value = -2347.625

$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 36000
$ws.Range("N125").Value = -40920
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23386.057
$ws.Range("I32").Value = 4478.493
$ws.Range("J32").Value = 102352.94
$ws.Range("K32").Value = 4478.493
$ws.Range("L32").Value = 102352.94
$ws.Range("M32").Value = -4191.493
$ws.Range("N32").Value = -102926.94

$ws.Range("H74").Value = 2481.4
$ws.Range("I74").Value = 1800
$ws.Range("K74").Value = 1800
$ws.Range("M74").Value = -926

$ws.Range("H77").Value = 2481.4
$ws.Range("I77").Value = 1800
$ws.Range("K77").Value = 9000
$ws.Range("M77").Value = -4632

$ws.Range("H97").Value = 671.6896400000001
$ws.Range("I97").Value = 703.5599999999999
$ws.Range("J97").Value = 472.5
$ws.Range("K97").Value = 703.5599999999999
$ws.Range("L97").Value = 472.5
$ws.Range("M97").Value = -207.5599999999999
$ws.Range("N97").Value = -1464.5

$ws.Range("H102").Value = 1303.8462
$ws.Range("I102").Value = 1204.1666
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1204.1666
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 417.8334
$ws.Range("N102").Value = -5744

$ws.Range("H122").Value = 2170.5
$ws.Range("I122").Value = 1116.6666
$ws.Range("J122").Value = 2802.8
$ws.Range("K122").Value = 3349.9998
$ws.Range("L122").Value = 8408.400000000001
$ws.Range("M122").Value = -899.9998000000001
$ws.Range("N122").Value = -13308.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2903.7097
$ws.Range("I86").Value = 2870.4565
$ws.Range("J86").Value = 2999.3125
$ws.Range("K86").Value = 2870.4565
$ws.Range("L86").Value = 2999.3125
$ws.Range("M86").Value = -1747.4565
$ws.Range("N86").Value = -5245.3125

$ws.Range("H89").Value = 2903.7097
$ws.Range("I89").Value = 2870.4565
$ws.Range("J89").Value = 2999.3125
$ws.Range("K89").Value = 14352.2825
$ws.Range("L89").Value = 14996.5625
$ws.Range("M89").Value = -8736.282499999999
$ws.Range("N89").Value = -26228.5625

$ws.Range("H94").Value = 872.0454999999999
$ws.Range("I94").Value = 808.75
$ws.Range("J94").Value = 1505
$ws.Range("K94").Value = 808.75
$ws.Range("L94").Value = 1505
$ws.Range("M94").Value = -357.75
$ws.Range("N94").Value = -2407

$ws.Range("H99").Value = 1374.5161
$ws.Range("I99").Value = 1250.4166
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1250.4166
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = 247.5834
$ws.Range("N99").Value = -4796

$ws.Range("H103").Value = 13326.333
$ws.Range("J103").Value = 13326.333
$ws.Range("L103").Value = 13326.333
$ws.Range("N103").Value = -15670.333

$ws.Range("H105").Value = 1570.375
$ws.Range("I105").Value = 1582.1305
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1582.1305
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 164.8695
$ws.Range("N105").Value = -4794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 862.63635
$ws.Range("I16").Value = 888.9
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 888.9
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -601.9
$ws.Range("N16").Value = -1174

$ws.Range("H105").Value = 746.6667
$ws.Range("I105").Value = 713.75
$ws.Range("J105").Value = 812.5
$ws.Range("K105").Value = 713.75
$ws.Range("L105").Value = 812.5
$ws.Range("M105").Value = 1033.25
$ws.Range("N105").Value = -4306.5

$ws.Range("H113").Value = 862.63635
$ws.Range("I113").Value = 888.9
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 888.9
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1281.1
$ws.Range("N113").Value = -4940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 495.6154
$ws.Range("I92").Value = 502
$ws.Range("J92").Value = 492.77777
$ws.Range("K92").Value = 1506
$ws.Range("L92").Value = 1478.33331
$ws.Range("M92").Value = -258
$ws.Range("N92").Value = -3974.33331

$ws.Range("H97").Value = 682.6875
$ws.Range("I97").Value = 212.3
$ws.Range("J97").Value = 1466.6666
$ws.Range("K97").Value = 636.9000000000001
$ws.Range("L97").Value = 4399.9998
$ws.Range("M97").Value = -140.9000000000001
$ws.Range("N97").Value = -5391.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2966
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2957.5
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2957.5
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4953.5

$ws.Range("H83").Value = 2966
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2957.5
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 14787.5
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24771.5

$ws.Range("H132").Value = 2219.468
$ws.Range("I132").Value = 2084.9697
$ws.Range("J132").Value = 2536.5
$ws.Range("K132").Value = 6254.909100000001
$ws.Range("L132").Value = 7609.5
$ws.Range("M132").Value = -3724.909100000001
$ws.Range("N132").Value = -12669.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2935
$ws.Range("I7").Value = 1800
$ws.Range("J7").Value = 3502.5
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 3502.5
$ws.Range("M7").Value = -1688
$ws.Range("N7").Value = -3726.5

$ws.Range("H61").Value = 2955
$ws.Range("J61").Value = 3057.8948
$ws.Range("L61").Value = 3057.8948
$ws.Range("N61").Value = -3461.8948

$ws.Range("H113").Value = 2955
$ws.Range("J113").Value = 3057.8948
$ws.Range("L113").Value = 3057.8948
$ws.Range("N113").Value = -7397.8948

$ws.Range("H126").Value = 2935
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 3502.5
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 10507.5
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -15447.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 345992.06
$ws.Range("I126").Value = 589294.1
$ws.Range("J126").Value = 1314.1666
$ws.Range("K126").Value = 1767882.3
$ws.Range("L126").Value = 3942.4998
$ws.Range("M126").Value = -1765412.3
$ws.Range("N126").Value = -8882.4998
